$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two extra data rows (rows 3 and 4), leaving header + one data row
$ws.Rows("3:4").Delete()

# Update row 2 with the new TPM-derived values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Guca2b"
$ws.Range("C2").Value = "Gucy2c"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2104683333333333
$ws.Range("H2").Value = 0.631405
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.015763
$ws.Range("N2").Value = 0.047289
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.003317612338333333
$ws.Range("R2").Value = 0.029858511045
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
